$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new data rows just before current row 692 (shifts rows 692..733 down to 694..735)
$ws.Rows.Item(692).Resize(2).Insert()

# Column A holds plain text dates ("2026/01/20" style), not real Excel dates.
# Force Text number format first so Value assignment doesn't auto-convert the
# string into a date serial, then restore the default ("Normal") style so the
# cell matches the rest of the sheet (no explicit style index).
$ws.Range("A692:A693").NumberFormat = "@"

# New row 692: 2026/01/20, 火, 22, 179
$ws.Cells.Item(692, 1).Value = "2026/01/20"
$ws.Cells.Item(692, 2).Value = "火"
$ws.Cells.Item(692, 3).Value = 22
$ws.Cells.Item(692, 4).Value = 179

# New row 693: 2026/01/21, 水, 2, 193
$ws.Cells.Item(693, 1).Value = "2026/01/21"
$ws.Cells.Item(693, 2).Value = "水"
$ws.Cells.Item(693, 3).Value = 2
$ws.Cells.Item(693, 4).Value = 193

$ws.Range("A692:A693").Style = "Normal"
